$d = $word.ActiveDocument

# --- Create the three new character styles -----------------------------
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Helper: apply a character style to every occurrence of some text --
function Apply-StyleToAllMatches($searchText, $style) {
    $pos = 0
    while ($true) {
        $range = $d.Range($pos, $d.Content.End)
        $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $range.Style = $style
        $pos = $range.End
    }
}

# 1) "Waktu Kampanye 2022 ..." runs -> GaNStyle (4 occurrences)
Apply-StyleToAllMatches "Waktu Kampanye 2022 untuk konstelasi sepatu bot: 14-23 Mei, 13-22 Juni, 12-21 Juli" $gaNStyle

# 2) "Anda sedang berpartisipasi ..." run -> GaNParagraph
Apply-StyleToAllMatches "Anda sedang berpartisipasi dalam kampanye global pengamatan dan pencatatan penampakan bintang paling redup untuk pengukuran tingkat polusi cahaya di suatu lokasi. Melalui pengamatan dan identifikasi  konstelasi sepatu bot di langit malam dan membandingkannya dengan peta bintang, masyarakat di seluruh dunia dapat mengetahui dan mempelajari seberapa besar kontribusi cahaya di lingkungannya terhadap polusi cahaya. Kontribusi data anda pada basis data online akan membantu mendokumentasikan langit malam yang tampak di berbagai lokasi." $gaNParagraph

# 3) "Peta di dokumen ini ..." run -> GaNLinks
Apply-StyleToAllMatches "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)." $gaNLinks
